$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a trailing-zero numeric string (e.g. "94.10").
# Excel infers Range.Value assignments that look numeric as numbers, which would
# silently drop the significant trailing zero (94.10 -> 94.1). Force these specific
# cells to Text number format first so the literal text is preserved, matching the
# source data (all Price/Volume cells in this sheet are plain text).
$textPriceCells = @("D6", "D26", "D27")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.602.36"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.235.84"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "269.65"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").Value = "94.10"
$ws.Range("E6").Value = "  +13.36%  "
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "46.45"
$ws.Range("E10").Value = "  +4.99%  "
$ws.Range("D11").Value = "0.0924"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "8.16"
$ws.Range("E12").Value = "  +15.19%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "2.572.55"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "15.16"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "2.232.53"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "0.804"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "43.592.96"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "70.42"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "2.33"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "233.24"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "2.50"
$ws.Range("E26").Value = "  +11.40%  "
$ws.Range("D27").Value = "11.20"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("E28").Value = "  +5.45%  "
$ws.Range("D29").Value = "40.46"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").Value = "172.83"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "0.0929"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").Value = "20.81"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("D37").Value = "0.0351"
$ws.Range("E37").Value = "  -4.35%  "
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  +19.63%  "
$ws.Range("D40").Value = "12.65"
$ws.Range("E40").Value = "  -6.13%  "
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").Value = "0.217"
$ws.Range("E42").Value = "  +7.19%  "
$ws.Range("D43").Value = "63.19"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "8.38"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0986"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "100.42"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "1.16"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.458.26"
$ws.Range("E51").Value = "  +0.14%  "
